$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: status columns (zh-cn / de-de) for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status
$wsZh.Range("K2").Value = "2016-10-26 08:33:12"
$wsZh.Range("K3").Value = "2016-10-26 08:33:12"
$wsZh.Range("P3").Value = ""

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status
$wsDe.Range("K2").Value = "2016-10-26 08:33:29"
$wsDe.Range("K3").Value = "2016-10-26 08:33:29"
$wsDe.Range("P3").Value = ""
